# Update "想去人数" (F column) figures across the four sheets to the
# newly scraped counts (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        3  = 821
        5  = 2283
        6  = 1360
        9  = 1143
        11 = 3001
        14 = 1099
        15 = 613
        17 = 230
        19 = 1068
        20 = 1068
        21 = 145
        22 = 524
        23 = 172
        25 = 209
        26 = 641
        28 = 3
        29 = 833
        30 = 66
        32 = 1035
        33 = 5029
        34 = 491
        35 = 230
        36 = 120
        37 = 3
        38 = 2
    }
    "演出" = @{
        3  = 20
        11 = 194
        20 = 40
        22 = 43
        27 = 656
        36 = 433
        37 = 433
        43 = 752
    }
    "本地生活" = @{
        4 = 638
        5 = 420
        6 = 401
    }
    "全部类型" = @{
        4  = 420
        5  = 821
        6  = 20
        9  = 2283
        10 = 1360
        16 = 3001
        19 = 1099
        20 = 613
        22 = 401
        26 = 230
        27 = 1068
        28 = 1068
        29 = 145
        31 = 524
        32 = 172
        33 = 209
        34 = 43
        35 = 641
        38 = 656
        39 = 833
        40 = 66
        41 = 1035
        42 = 5029
        44 = 491
        46 = 433
        47 = 230
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
